$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6041666666666666
$ws.Range("D2").Value = 0.6744186046511628
$ws.Range("B3").Value = 0.8
$ws.Range("C3").Value = 0.6545454545454545
$ws.Range("D3").Value = 0.7200000000000001
$ws.Range("B4").Value = 0.6989247311827957
$ws.Range("C4").Value = 0.6989247311827957
$ws.Range("D4").Value = 0.6989247311827957
$ws.Range("E4").Value = 0.6989247311827957
$ws.Range("B5").Value = 0.7020833333333334
$ws.Range("C5").Value = 0.7088516746411484
$ws.Range("D5").Value = 0.6972093023255814
$ws.Range("B6").Value = 0.7199820788530465
$ws.Range("C6").Value = 0.6989247311827957
$ws.Range("D6").Value = 0.7013753438359589
$ws.Range("B12").Value = 0.5510204081632653
$ws.Range("C12").Value = 0.7105263157894737
$ws.Range("D12").Value = 0.6206896551724138
$ws.Range("B13").Value = 0.75
$ws.Range("C13").Value = 0.6
$ws.Range("D13").Value = 0.6666666666666665
$ws.Range("B14").Value = 0.6451612903225806
$ws.Range("C14").Value = 0.6451612903225806
$ws.Range("D14").Value = 0.6451612903225806
$ws.Range("E14").Value = 0.6451612903225806
$ws.Range("B15").Value = 0.6505102040816326
$ws.Range("C15").Value = 0.6552631578947368
$ws.Range("D15").Value = 0.6436781609195401
$ws.Range("B16").Value = 0.6686965108624094
$ws.Range("C16").Value = 0.6451612903225806
$ws.Range("D16").Value = 0.6478803608948214
$ws.Range("B17").Value = 0.6326530612244898
$ws.Range("C17").Value = 0.8157894736842105
$ws.Range("D17").Value = 0.7126436781609196
$ws.Range("B18").Value = 0.8409090909090909
$ws.Range("C18").Value = 0.6727272727272727
$ws.Range("D18").Value = 0.7474747474747475
$ws.Range("B19").Value = 0.7311827956989247
$ws.Range("C19").Value = 0.7311827956989247
$ws.Range("D19").Value = 0.7311827956989247
$ws.Range("E19").Value = 0.7311827956989247
$ws.Range("B20").Value = 0.7367810760667903
$ws.Range("C20").Value = 0.7442583732057416
$ws.Range("D20").Value = 0.7300592128178336
$ws.Range("B21").Value = 0.7558152293175335
$ws.Range("C21").Value = 0.7311827956989247
$ws.Range("D21").Value = 0.733242697647592
$ws.Range("B22").Value = 0.6666666666666666
$ws.Range("C22").Value = 0.7894736842105263
$ws.Range("D22").Value = 0.7228915662650601
$ws.Range("B23").Value = 0.8333333333333334
$ws.Range("C23").Value = 0.7272727272727273
$ws.Range("D23").Value = 0.7766990291262137
$ws.Range("B24").Value = 0.7526881720430108
$ws.Range("C24").Value = 0.7526881720430108
$ws.Range("D24").Value = 0.7526881720430108
$ws.Range("E24").Value = 0.7526881720430108
$ws.Range("B25").Value = 0.75
$ws.Range("C25").Value = 0.7583732057416268
$ws.Range("D25").Value = 0.7497952976956369
$ws.Range("B26").Value = 0.7652329749103943
$ws.Range("C26").Value = 0.7526881720430108
$ws.Range("D26").Value = 0.7547131840861725
